# Generate Report for Handback
#
# Applies the localization-status handback update:
#  - "Ready for handoff" status text becomes "Handed back: in sync with en-US"
#    everywhere it is used (Overview + per-language sheets).
#  - The zh-cn sheet gets its handback file / handback datetime populated,
#    including a hyperlink to the source doc in the new "Latest Target File" column.
#  - The de-de sheet gets the same treatment with its own handback datetime.
#  - A handful of columns are widened so the new, longer text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Update every cell currently carrying that status so the change is global.
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# --- zh-cn: populate Latest Target File (I), Latest Handback File (J),
#     and Latest Handback DateTime (K) ---
$zhCnDocUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0342fff73d79c6b90cae3bb4114e85b51ded9781/e2e/714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhCnDocUrl, "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhCnDocUrl, "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md")

$wsZhCn.Range("J2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-09-03 07:06:43"
$wsZhCn.Range("K3").Value = "2016-09-03 07:06:43"

# --- de-de: same idea, with its own xlf file name and handback datetime ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $zhCnDocUrl, "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $zhCnDocUrl, "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md")

$wsDeDe.Range("J2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.de-de.xlf"
$wsDeDe.Range("J3").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-09-03 07:06:53"
$wsDeDe.Range("K3").Value = "2016-09-03 07:06:53"

# --- Column widths: widen columns that now hold the longer text ---
# (ColumnWidth is specified in characters; the host snaps to the nearest
# pixel, so the inputs below are chosen to land on the intended stored width.)
$wideStatusWidth = 29.166666666666668   # -> stored width ~29.98
$wideFileWidth   = 39.166666666666664   # -> stored width 40

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth  = $wideStatusWidth   # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth  = $wideFileWidth     # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = $wideFileWidth     # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth  = $wideStatusWidth   # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth  = $wideFileWidth     # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = $wideFileWidth     # J: Latest Handback File

Write-Host "Handback report generated."
